$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 380, shifting rows 380:469 down to 381:470
$ws.Rows("380:380").Insert()

# Populate the new row 380 with the new data record
$ws.Range("A380").Value = 11
$ws.Range("B380").Value = "Vega Monumental Concepción"
$ws.Range("C380").Value = "Bíobío"
$ws.Range("D380").Value = 44642
$ws.Range("D380").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E380").Value = 8
$ws.Range("F380").Value = "Fruta"
$ws.Range("G380").Value = 100108
$ws.Range("H380").Value = "Tropicales y subtropicales"
$ws.Range("I380").Value = 100108006
$ws.Range("J380").Value = "Plátano"
$ws.Range("K380").Value = "Sin especificar"
$ws.Range("L380").Value = "Pintón"
$ws.Range("M380").Value = 1100
$ws.Range("N380").Value = 18000
$ws.Range("O380").Value = 19000
$ws.Range("P380").Value = 18455
$ws.Range("Q380").Value = "$/caja 20 kilos"
$ws.Range("R380").Value = "Ecuador"
$ws.Range("S380").Value = 923
$ws.Range("T380").Value = 20
